$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = -5
$ws.Range("F11").Value = -7
$ws.Range("F19").Value = 8
$ws.Range("F20").Value = 10
$ws.Range("F21").Value = -3
